$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3493
$ws1.Range("F3").Value = 745
$ws1.Range("F5").Value = 7012
$ws1.Range("F6").Value = 2787
$ws1.Range("F7").Value = 52
$ws1.Range("F8").Value = 125
$ws1.Range("F11").Value = 86
$ws1.Range("F12").Value = 44
$ws1.Range("F13").Value = 4
$ws1.Range("F14").Value = 178
$ws1.Range("F16").Value = 7

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("C2").Value = "合肥·全国地下偶像联合公演展-永乐大典Vol.01（取消）"
$ws2.Range("F2").Value = 28
$ws2.Range("G2").Value = "不可售"

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3493
$ws4.Range("C3").Value = "合肥·全国地下偶像联合公演展-永乐大典Vol.01（取消）"
$ws4.Range("F3").Value = 28
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F4").Value = 745
$ws4.Range("F6").Value = 7012
$ws4.Range("F7").Value = 2787
$ws4.Range("F8").Value = 52
$ws4.Range("F9").Value = 125
$ws4.Range("F12").Value = 86
$ws4.Range("F13").Value = 44
$ws4.Range("F14").Value = 4
$ws4.Range("F15").Value = 178
$ws4.Range("F17").Value = 7
